$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Fenêtre de jeu et chonomètre :" -> "Fenêtre de jeu et minuteur :"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Fenêtre de jeu et chonomètre :", $true, $true, $false, $false, $false,
    $true, 1, $false, "Fenêtre de jeu et minuteur :", 2)

# ---------------------------------------------------------------------------
# 2) Body paragraph right after: drop the stray "le " before "des canvas"
#    and rename "chronomètre" -> "minuteur"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "avec le des canvas et des lines", $true, $true, $false, $false, $false,
    $true, 1, $false, "avec des canvas et des lines", 2)

$null = $d.Content.Find.Execute(
    "j'ai ajouté le chronomètre en haut à droite",
    $true, $true, $false, $false, $false,
    $true, 1, $false, "j'ai ajouté le minuteur en haut à droite", 2)

# ---------------------------------------------------------------------------
# 3) "Pour les fenêtres de fin je l'ai aussi créé ..." paragraph
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Pour les fenêtres de fin je l'ai aussi créé",
    $true, $true, $false, $false, $false,
    $true, 1, $false, "Pour les fenêtres de fin je les ai aussi créé", 2)

$null = $d.Content.Find.Execute(
    "le chronomètre est à 0 ou lorsque tous les objets",
    $true, $true, $false, $false, $false,
    $true, 1, $false, "le chronomètre est à 0 ou que tous les objets", 2)

# ---------------------------------------------------------------------------
# 4) "Emplacement des items ..." paragraph - full rewrite
# ---------------------------------------------------------------------------
$oldItems = "Emplacement des items aléatoires avec certaines distances en fonctions de la difficulté. Ensuite affichage sur la minimap en rouge. Emplacement de l'antenne dans l'emplacement de la difficulté. Changement de skin en fonction de l'avancé des réparations de l'antenne et affichage en bleu sur la minimap."
$newItems = "Emplacement des items aléatoires avec certaines distances en fonction de la difficulté. Ensuite affichage sur la minimap en rouge (affichage ou non des objets en fonction s'ils ont déjà servi à réparer l'antenne ou s'ils sont dans l'inventaire). Emplacement de l'antenne en fonction de la difficulté. Changement de skin de l'antenne en fonction de l'avancé des réparations et affichage en bleu sur la minimap."
$null = $d.Content.Find.Execute(
    $oldItems, $true, $true, $false, $false, $false,
    $true, 1, $false, $newItems, 2)

# ---------------------------------------------------------------------------
# 5) New sections appended after the "Mise en place des items sur la map"
#    paragraph (the one we just rewrote above).
# ---------------------------------------------------------------------------

# Find the paragraph that now ends with "...sur la minimap." (the rewritten one)
$lastPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*affichage en bleu sur la minimap.*") {
        $lastPara = $para
    }
}

$anchorIndex = $lastPara.Index

# -- "Drag and drop :  " (Titre2) ------------------------------------------
$null = $d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($anchorIndex + 1)
$p1.Style = "Titre2"
$p1.Range.InsertBefore("Drag and drop :  ")

# -- Drag and drop body text (Normal) ---------------------------------------
$null = $d.Paragraphs($anchorIndex + 1).Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($anchorIndex + 2)
$p2.Style = "Normal"
$p2.Range.InsertBefore("Pour le drag and drop je vérifie si le clic est fait sur un objet et ensuite s'il est fait sur un objet alors je récupère ses coordonnées dynamiquement pour qu'il s'affiche à l'endroit où est la souris puis lors du release je vérifie si l'objet est lâché sur l'antenne (si oui je fais disparaitre l'objet et je passe à l'étape de réparation d'après, sinon il est toujours visible). ")

# -- "Inventaire : " (Titre2) -------------------------------------------------
$null = $d.Paragraphs($anchorIndex + 2).Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($anchorIndex + 3)
$p3.Style = "Titre2"
$p3.Range.InsertBefore("Inventaire : ")

# -- Inventaire body text (Normal) -------------------------------------------
$null = $d.Paragraphs($anchorIndex + 3).Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($anchorIndex + 4)
$p4.Style = "Normal"
$p4.Range.InsertBefore("J'ai eu du mal à le faire car je ne savais pas comment j'allais sortir les objets de l'inventaire. J'ai donc décidé de faire des canvas et lorsque je clique sur un canvas et que l'objet est dans cet emplacement alors il est drop près de la case. Ensuite je peux faire le drag and drop.")

# -- "Mouvement du véhicule :" (Titre2) --------------------------------------
$null = $d.Paragraphs($anchorIndex + 4).Range.InsertParagraphAfter()
$p5 = $d.Paragraphs($anchorIndex + 5)
$p5.Style = "Titre2"
$p5.Range.InsertBefore("Mouvement du véhicule :")

# -- Mouvement body text (Normal) --------------------------------------------
$null = $d.Paragraphs($anchorIndex + 5).Range.InsertParagraphAfter()
$p6 = $d.Paragraphs($anchorIndex + 6)
$p6.Style = "Normal"
$p6.Range.InsertBefore("Je récupère sur quelle flèche on clique pour orienter le véhicule dans la bonne direction.")

# -- trailing empty paragraph -------------------------------------------------
$null = $d.Paragraphs($anchorIndex + 6).Range.InsertParagraphAfter()
$p7 = $d.Paragraphs($anchorIndex + 7)
$p7.Style = "Normal"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
